# Update the lattice-multiplication exercise table: every cell's three
# pieces of text (problem header, second-factor digits, first-factor
# digits running down the left edge) are replaced with a new problem,
# while the table's shape (5 rows x 3 columns) and each run's formatting
# (sz=32) stay untouched.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$br = [char]11

# New content per cell, addressed (row, col) 1-based, matching Word's
# Table.Cell(r, c) convention. Each tuple is:
#   header, second-factor row, left-digit-1 row, left-digit-2 row
$cells = @(
    @{ r = 1; c = 1; top = "49 x 50"; mid = "  5    0"; l1 = "4|    |"; l2 = "9|    |" },
    @{ r = 1; c = 2; top = "68 x 21"; mid = "  2    1"; l1 = "6|    |"; l2 = "8|    |" },
    @{ r = 1; c = 3; top = "43 x 95"; mid = "  9    5"; l1 = "4|    |"; l2 = "3|    |" },

    @{ r = 2; c = 1; top = "70 x 57"; mid = "  5    7"; l1 = "7|    |"; l2 = "0|    |" },
    @{ r = 2; c = 2; top = "61 x 50"; mid = "  5    0"; l1 = "6|    |"; l2 = "1|    |" },
    @{ r = 2; c = 3; top = "25 x 14"; mid = "  1    4"; l1 = "2|    |"; l2 = "5|    |" },

    @{ r = 3; c = 1; top = "93 x 32"; mid = "  3    2"; l1 = "9|    |"; l2 = "3|    |" },
    @{ r = 3; c = 2; top = "38 x 96"; mid = "  9    6"; l1 = "3|    |"; l2 = "8|    |" },
    @{ r = 3; c = 3; top = "89 x 21"; mid = "  2    1"; l1 = "8|    |"; l2 = "9|    |" },

    @{ r = 4; c = 1; top = "92 x 95"; mid = "  9    5"; l1 = "9|    |"; l2 = "2|    |" },
    @{ r = 4; c = 2; top = "70 x 49"; mid = "  4    9"; l1 = "7|    |"; l2 = "0|    |" },
    @{ r = 4; c = 3; top = "66 x 65"; mid = "  6    5"; l1 = "6|    |"; l2 = "6|    |" },

    @{ r = 5; c = 1; top = "34 x 61"; mid = "  6    1"; l1 = "3|    |"; l2 = "4|    |" },
    @{ r = 5; c = 2; top = "20 x 54"; mid = "  5    4"; l1 = "2|    |"; l2 = "0|    |" },
    @{ r = 5; c = 3; top = "57 x 50"; mid = "  5    0"; l1 = "5|    |"; l2 = "7|    |" }
)

foreach ($cell in $cells) {
    $tc = $t.Cell($cell.r, $cell.c)
    $newText = $cell.top + $br + $cell.mid + $br + "  ----" + $br + $cell.l1 + $br + $cell.l2
    $tc.Range.Text = $newText
}
